# Update projection results on Sheet1 (rows 2-6) to reflect the revised
# headcount / contribution projections.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 88
$ws.Range("E2").Value = 0.8627450980392157
$ws.Range("F2").Value = 0.8627450980392157
$ws.Range("G2").Value = 0.09737246067906356
$ws.Range("H2").Value = 0.08400761313487837
$ws.Range("I2").Value = 462691.8401053585
$ws.Range("J2").Value = 168731.9202536792
$ws.Range("L2").Value = 168731.9202536792
$ws.Range("M2").Value = 631423.7603590377
$ws.Range("N2").Value = 10131360.5088
$ws.Range("O2").Value = 9723619.5787
$ws.Range("P2").Value = 0.01665441873350774
$ws.Range("Q2").Value = 0.01735278914276877

# Row 3
$ws.Range("D3").Value = 88
$ws.Range("E3").Value = 0.8543689320388349
$ws.Range("F3").Value = 0.8543689320388349
$ws.Range("G3").Value = 0.09684574533725648
$ws.Range("H3").Value = 0.08274199601629682
$ws.Range("I3").Value = 477113.6797688863
$ws.Range("J3").Value = 173344.4200914731
$ws.Range("L3").Value = 173344.4200914731
$ws.Range("M3").Value = 650458.0998603592
$ws.Range("N3").Value = 10494911.028964
$ws.Range("O3").Value = 10087537.870961
$ws.Range("P3").Value = 0.0165169975822639
$ws.Range("Q3").Value = 0.01718401678475774

# Row 4
$ws.Range("D4").Value = 88
$ws.Range("E4").Value = 0.8461538461538461
$ws.Range("F4").Value = 0.8461538461538461
$ws.Range("G4").Value = 0.0971597766805547
$ws.Range("H4").Value = 0.08221211872970012
$ws.Range("I4").Value = 504436.7308120827
$ws.Range("J4").Value = 179919.5730192822
$ws.Range("L4").Value = 179919.5730192822
$ws.Range("M4").Value = 684356.3038313651
$ws.Range("N4").Value = 10907997.05983292
$ws.Range("O4").Value = 10499652.70708983
$ws.Range("P4").Value = 0.01649428140036903
$ws.Range("Q4").Value = 0.01713576420463818

# Row 5
$ws.Range("D5").Value = 89
$ws.Range("E5").Value = 0.8476190476190476
$ws.Range("F5").Value = 0.8476190476190476
$ws.Range("G5").Value = 0.0964615420761171
$ws.Range("H5").Value = 0.08176264042642307
$ws.Range("I5").Value = 524764.3891885336
$ws.Range("J5").Value = 187509.438435905
$ws.Range("L5").Value = 187509.438435905
$ws.Range("M5").Value = 712273.8276244387
$ws.Range("N5").Value = 11375845.10502791
$ws.Range("O5").Value = 10965150.42170252
$ws.Range("P5").Value = 0.01648312162346773
$ws.Range("Q5").Value = 0.01710048938907224

# Row 6
$ws.Range("D6").Value = 91
$ws.Range("E6").Value = 0.8584905660377359
$ws.Range("F6").Value = 0.8584905660377359
$ws.Range("G6").Value = 0.0955045359716872
$ws.Range("H6").Value = 0.08198974314550507
$ws.Range("I6").Value = 548908.0304184185
$ws.Range("J6").Value = 196520.0763660966
$ws.Range("L6").Value = 196520.0763660966
$ws.Range("M6").Value = 745428.1067845151
$ws.Range("N6").Value = 11767757.80777875
$ws.Range("O6").Value = 11353292.2839536
$ws.Range("P6").Value = 0.01669987431558053
$ws.Range("Q6").Value = 0.0173095232159091
